# The sheet currently has only a header row (A1:F1). The edit extends it
# with two more "possible value" template rows (A1:F3 overall). Columns
# B (light_route_type) and F (light_stick_type) get the sample values;
# the other columns (A, C, D, E) stay blank placeholders but are still
# materialized as real (empty) cells, matching the template's row shape.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: train / mast
$ws.Range("A2").Style = "Normal"
$ws.Range("B2").Value = "train"
$ws.Range("C2:E2").Style = "Normal"
$ws.Range("F2").Value = "mast"

# Row 3: shunt / dwarf
$ws.Range("A3").Style = "Normal"
$ws.Range("B3").Value = "shunt"
$ws.Range("C3:E3").Style = "Normal"
$ws.Range("F3").Value = "dwarf"
